$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 25: "GFG" | "GFG" | "Diagonal Traversal of Binary Tree "
$ws.Cells.Item(25, 1).Value = "GFG"
$ws.Cells.Item(25, 2).Value = "GFG"
$ws.Cells.Item(25, 3).Value = "Diagonal Traversal of Binary Tree "

# Match the formatting used by the row above (row 24), which has the same
# GFG / GFG / <question text> pattern.
$ws.Cells.Item(25, 1).HorizontalAlignment = -4131
$ws.Cells.Item(25, 1).VerticalAlignment = -4107
$ws.Cells.Item(25, 1).WrapText = $False

$ws.Cells.Item(25, 2).HorizontalAlignment = -4131
$ws.Cells.Item(25, 2).VerticalAlignment = -4160
$ws.Cells.Item(25, 2).WrapText = $True

$ws.Cells.Item(25, 3).HorizontalAlignment = -4131
$ws.Cells.Item(25, 3).VerticalAlignment = -4160
$ws.Cells.Item(25, 3).WrapText = $True

# Update the saved selection to C27 (was D27)
$null = $ws.Range("C27").Select()
